$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Turn the "Deployed Website Link: <hyperlink>" paragraph into a
# plain-text "Repository link: <url>" paragraph (no hyperlink anymore).
# ---------------------------------------------------------------------------

# First change the display/target text of the hyperlink run itself while it
# is still a real hyperlink field (this also strips the Hyperlink character
# style from the run, same as Word does when the link text is edited).
$hyperlink = $d.Hyperlinks(1)
$linkRange = $hyperlink.Range
$null = $linkRange.Find.Execute(
    "https://tjc-ltd.github.io/Wheels-On-Fire-Website-Prototype/html/index.html",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://github.com/TJC-Ltd/Wheels-On-Fire-Website-Prototype", 2)

# Update the leading label text in the first run of the paragraph.
$null = $d.Content.Find.Execute(
    "Deployed Website Link: ", $true, $false, $false, $false, $false, $true,
    1, $false, "Repository link: ", 2)

# Finally, remove the hyperlink field itself, leaving the (already updated)
# display text behind as plain text runs.
$hyperlink2 = $d.Hyperlinks(1)
$hyperlink2.Delete()

# ---------------------------------------------------------------------------
# Edit 2: Break "in order to" out of the "Secondarily though, ..." sentence
# into its own run (mirrors the grammar-checker gramStart/gramEnd wrapping
# that Word adds around that phrase - the visible text is unchanged).
# ---------------------------------------------------------------------------

$finder = $d.Content
$null = $finder.Find.Execute(
    "in order to inconvenience", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$phraseStart = $finder.Start

$phrase = $d.Range($phraseStart, $phraseStart + 11)

# Toggling a character-formatting property and then reverting it forces Word
# to split this sub-range into its own run, without altering the visible
# text or the run's original rsid.
$phrase.Font.Bold = 1
$phrase.Font.Bold = 0
